$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 213; existing rows 213:242 shift down to 214:243.
$ws.Rows.Item(213).Insert()

# Populate the newly inserted row 213 with the new weekly record
# (same market/category/etc. as the surrounding rows, new date + prices).
$ws.Cells.Item(213, 1).Value = 8
$ws.Cells.Item(213, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(213, 3).Value = "Coquimbo"
$ws.Cells.Item(213, 4).Value = [DateTime]"2023-07-20"
$ws.Cells.Item(213, 5).Value = 4
$ws.Cells.Item(213, 6).Value = 100112044
$ws.Cells.Item(213, 7).Value = "Perejil"
$ws.Cells.Item(213, 8).Value = "Sin especificar"
$ws.Cells.Item(213, 9).Value = "Primera"
$ws.Cells.Item(213, 10).Value = 2400
$ws.Cells.Item(213, 11).Value = 2500
$ws.Cells.Item(213, 12).Value = 3000
$ws.Cells.Item(213, 13).Value = 2750
$ws.Cells.Item(213, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(213, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(213, 16).Value = 1833
$ws.Cells.Item(213, 17).Value = 1.5
$ws.Cells.Item(213, 18).Value = "Hortaliza"
